$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "Alexis Sharon_20251130_195228"
# B14: left blank (empty cell in source)
$ws.Range("C14").Value = "Alexis Sharon"
$ws.Range("D14").Value = 28
$ws.Range("E14").Value = "Female"
$ws.Range("F14").Value = "2025-11-30 19:52:28"
$ws.Range("G14").Value = "{
  `"portion`": 0.4,
  `"diet`": 0.2857142857142857,
  `"salt`": 0.6,
  `"fat`": 0.2,
  `"natural`": 0.4,
  `"convenience`": 0.2,
  `"price`": 0.6
}"
$ws.Range("H14").Value = "Nongshim Neoguri Spicy Seafood"
$ws.Range("I14").Value = "'0.566"
$ws.Range("J14").Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"
$ws.Range("K14").Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Range("L14").Value = "'0.555"
$ws.Range("M14").Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"
$ws.Range("N14").Value = "Maruchan Ramen Sabor Pollo"
$ws.Range("O14").Value = "'0.530"
$ws.Range("P14").Value = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"
$ws.Range("Q14").Value = "Kraft Macaroni & Cheese Dinner"
$ws.Range("R14").Value = "'0.674"
$ws.Range("S14").Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"
$ws.Range("T14").Value = "Annie’s Shells & White Cheddar"
$ws.Range("U14").Value = "'0.602"
$ws.Range("V14").Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"
$ws.Range("W14").Value = "Amy’s Macaroni & Cheese (frozen)"
$ws.Range("X14").Value = "'0.595"
$ws.Range("Y14").Value = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"
$ws.Range("Z14").Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Range("AA14").Value = "'0.637"
$ws.Range("AB14").Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"
$ws.Range("AC14").Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Range("AD14").Value = "'0.595"
$ws.Range("AE14").Value = "Portátil, saludable, fácil, buena textura, sabor suave"
$ws.Range("AF14").Value = "Jack Link’s Beef Jerky Original"
$ws.Range("AG14").Value = "'0.594"
$ws.Range("AH14").Value = "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"

$ws.Rows.Item(14).AutoFit()
